$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("数组")

$d8Text = @'
1 数组已经排序，设置前后指针low=0，high=nums.length-1
2 计算两个指针之和是否是目标值，如果是则返回
3 如果和大于target，说明high左移可以让和变小（数组有序），high左移
4 如果和小于target，说明low右移可以让和变大（数组有序），low右移
5 循环存在条件是low < high
6 返回low+1,high+1
'@
$c8Text = @'
给定一个已按照升序排列 的有序数组，找到两个数使得它们相加之和等于目标数。 
 函数应该返回这两个下标值 index1 和 index2，其中 index1 必须小于 index2。 
 说明: 
 返回的下标值（index1 和 index2）不是从零开始的。 
 你可以假设每个输入只对应唯一的答案，而且你不可以重复使用相同的元素。 
 示例: 
 输入: numbers = [2, 7, 11, 15], target = 9
输出: [1,2]
解释: 2 与 7 之和等于目标数 9 。因此 index1 = 1, index2 = 2 。 
 Related Topics 数组 双指针 二分查找
'@
$e8Text = @'
排序
左右双指针
'@

$ws.Range("A8").Value = 7
$ws.Range("B8").Value = 167
$ws.Range("D8").Value = $d8Text
$ws.Range("C8").Value = $c8Text
$ws.Range("E8").Value = $e8Text
$ws.Range("F8").Value = "O(N)"
$ws.Range("G8").Value = "O(1)"

$ws.Rows.Item(8).RowHeight = 374

$ws.Range("D8").Select()
